# Update BunkerPrices at 2025-04-18 06:01
# Append a new data row (row 26) to the bunker prices table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 26
$prevRow = $newRow - 1

# Date value for column A (serial date number), formatted like the row above it.
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($prevRow, 1).NumberFormat
$ws.Cells.Item($newRow, 1).Value = 45762

# Bunker price values for columns B (2) through AU (47).
$values = @(485, 480, 520, 540, 535, 540, 505, 486, 485, 711, 721, 480, 480, 558, 560, 483, 625, 483, 490, 535, 445, 445, 477, 478, 570, 490, 490, 485, 560, 470, 492, 548, 494, 494, 527, 523.5, 546, 544, 502, 813, 600, 612, 611, 567, 594, 636)

$col = 2
foreach ($v in $values) {
    $ws.Cells.Item($newRow, $col).Value = $v
    $col = $col + 1
}
